$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above so the new row's date cell reuses
# the existing date-formatted style (numFmtId 22) instead of creating a
# new style entry.
$ws.Range("A3:N3").Copy($ws.Range("A4:N4"))

$ws.Range("A4").Value = 42606.881157407406
$ws.Range("B4").Value = 46
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = "Random"
